$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record arrived for "Repollo" (Crespo record) at the
# "Feria Lagunitas de Puerto Montt" market. Insert a fresh row at 454,
# pushing the existing rows 454:473 down to 455:474 (dimension becomes
# A1:R474), then populate the new row with its data.
$ws.Rows("454:454").Insert()

$ws.Cells.Item(454, 1).Value = 4
$ws.Cells.Item(454, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(454, 3).Value = "Los Lagos"
$ws.Cells.Item(454, 4).Value = 44753
$ws.Cells.Item(454, 5).Value = 10
$ws.Cells.Item(454, 6).Value = 100112006
$ws.Cells.Item(454, 7).Value = "Repollo"
$ws.Cells.Item(454, 8).Value = "Crespo record"
$ws.Cells.Item(454, 9).Value = "Primera"
$ws.Cells.Item(454, 10).Value = 500
$ws.Cells.Item(454, 11).Value = 2000
$ws.Cells.Item(454, 12).Value = 2000
$ws.Cells.Item(454, 13).Value = 2000
$ws.Cells.Item(454, 14).Value = "$/unidad"
$ws.Cells.Item(454, 15).Value = "Región del Maule"
$ws.Cells.Item(454, 16).Value = 2000
$ws.Cells.Item(454, 17).Value = 1
$ws.Cells.Item(454, 18).Value = "Hortaliza"
